$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.85"
$ws.Range("D3").Value = "'23.07"
$ws.Range("D4").Value = "'5.399"
$ws.Range("D5").Value = "'0.05974"
$ws.Range("D6").Value = "'3.403"
$ws.Range("D7").Value = "'6.481"
$ws.Range("D8").Value = "'0.8133"
$ws.Range("D9").Value = "'0.9071"
$ws.Range("D10").Value = "'0.1418"
$ws.Range("D11").Value = "'0.07422"
$ws.Range("D12").Value = "'0.03331"
$ws.Range("D14").Value = "'0.09342"
$ws.Range("D15").Value = "'3.855"
$ws.Range("D16").Value = "'0.001576"
$ws.Range("D17").Value = "'0.04634"
$ws.Range("D18").Value = "'0.0005940"
$ws.Range("D19").Value = "'0.006113"
$ws.Range("D20").Value = "'0.005027"
$ws.Range("E20").Value = "19HotbitTokenHTBBestin24h"
$ws.Range("D21").Value = "'0.0009816"
$ws.Range("D22").Value = "'0.00007796"
$ws.Range("D23").Value = "'0.0002900"
$ws.Range("D24").Value = "'3.621"
$ws.Range("D40").Value = "'0.03889"
$ws.Range("D41").Value = "'0.006200"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D43").Value = "'0.002569"
$ws.Range("D44").Value = "'0.007208"
$ws.Range("D45").Value = "'0.00005183"
$ws.Range("D47").Value = "'0.0005800"
$ws.Range("D48").Value = "'1.045"
$ws.Range("D49").Value = "'0.002261"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D51").Value = "'0.0002000"
